# Delete the empty placeholder slide (blank "ctrTitle" / "subTitle" / slide
# number layout) that currently sits at position 15, right before the
# "First rule of hooks / Callback and Ref / Custom hooks / Context" slide.
# Removing it shifts the latter slide up to become the new (last) slide 15,
# matching the "Added Hands On Demos - Day 18" commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$s.Delete()
